# REFACTOR: Bring all of Dokeza up to date including the online version.
#
# The document was round-tripped through an external converter/refactor
# that strips the SharePoint document-library "custom XML parts" that had
# been injected into the package (contentTypeSchema, the SharePoint
# FormTemplates part, and the empty documentManagement properties part).
# None of this metadata is visible document content — it only lives in
# customXml/item*.xml (+ their itemProps*.xml companions) — so the fix is
# to drop those custom XML parts from ActiveDocument entirely.

$d = $word.ActiveDocument

# Custom XML parts are attached to the document (not the visible body),
# so walk the collection back-to-front and delete every non-built-in part.
$customParts = $d.CustomXMLParts
for ($i = $customParts.Count; $i -ge 1; $i--) {
    $part = $customParts.Item($i)
    $part.Delete()
}

# Also sweep the "include built-in" view in case any of the injected
# SharePoint parts were (mis)classified as built-in by the host — only
# the three SharePoint-authored schemas/templates should ever match here;
# genuine Word built-ins (core/extended/cover-page properties) have none
# of these namespaces and are left untouched.
$sharePointNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms",
    "http://schemas.microsoft.com/office/2006/metadata/properties"
)

$allParts = $d.CustomXMLParts($true)
for ($i = $allParts.Count; $i -ge 1; $i--) {
    $part = $allParts.Item($i)
    if ($sharePointNamespaces -contains $part.NamespaceURI) {
        $part.Delete()
    }
}

Write-Output "CustomXMLParts remaining: $($d.CustomXMLParts($true).Count)"
